$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 97.67168086543172
$ws.Range("H2").Value = 97.43673423899337
$ws.Range("I2").Value = 95.99104519356429

$ws.Range("G3").Value = 97.98257399853244
$ws.Range("H3").Value = 97.28207661140405
$ws.Range("I3").Value = 96.19876235274121

$ws.Range("G4").Value = 97.82840967211807
$ws.Range("H4").Value = 97.26767829139126
$ws.Range("I4").Value = 95.84247630362455

$ws.Range("G5").Value = 97.74999335750779
$ws.Range("H5").Value = 97.20078650004831
$ws.Range("I5").Value = 96.00854362975288

$ws.Range("G6").Value = 97.86485448710926
$ws.Range("H6").Value = 97.29865438527905
$ws.Range("I6").Value = 95.89948081530729
